# ======================================================================
# Commit: "Added all Summary-scenarios except 5"
#
# - Renames RESUMEN E1's scenario caption from the leftover English
#   "Scenario 1: ..." text to the Spanish "Escenario: ..." wording used
#   by the rest of the workbook.
# - Duplicates RESUMEN E1 three times (so every formatting detail - the
#   bold/bordered/centered header row, column widths, etc. - matches
#   exactly) to create RESUMEN E2, RESUMEN E3 and RESUMEN E4, then
#   overwrites each duplicate's data rows with its own scenario figures.
# ======================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RESUMEN E1")

# --- Fix the RESUMEN E1 scenario label ("Scenario 1" -> "Escenario") ---
$ws1.Cells.Item(2, 1).Value2 = "Escenario: medidores con VAN > 0"

# --- Scenario data for RESUMEN E2 (rows 3-10, columns B-P) ---
$data5 = @{}
$data5[3] = @(1032, 1032, 1032, 1032, 1032, 1032, 1032, 1032, 1032, 1032, 1032, 1032, 1032, 1032, 1032)
$data5[4] = @(21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592, 21.06122448979592)
$data5[5] = @(62.83630668747166, 65.05964390846269, 66.94964062577307, 68.59710686745984, 70.0573338675697, 71.37202353252101, 72.56166063718823, 73.65058967517743, 74.6545472440211, 75.5883067854389, 76.14252432532948, 76.16546260021182, 76.16546260021182, 76.16546260021182, 76.16546260021182)
$data5[6] = @(41.77508219767574, 43.99841941866677, 45.88841613597715, 47.53588237766392, 48.99610937777378, 50.31079904272509, 51.50043614739231, 52.58936518538151, 53.59332275422518, 54.52708229564298, 55.08129983553356, 55.10423811041591, 55.10423811041591, 55.10423811041591, 55.10423811041591)
$data5[7] = @(0.04047973081170129, 0.04263412734366935, 0.04446551951160577, 0.04606190152874411, 0.04747685017226141, 0.04875077426620648, 0.04990352339863596, 0.05095868719513712, 0.05193151429672983, 0.05283632005391762, 0.05337335255381159, 0.0533955795643565, 0.0533955795643565, 0.0533955795643565, 0.0533955795643565)
$data5[8] = @(73515.3303129729, 77427.93470773762, 80753.93014937529, 83653.12312675768, 86222.81454560472, 88536.39096229758, 90629.90125841962, 92546.18660613867, 94312.94009676037, 95956.15987058138, 96931.46580337416, 96971.83232729399, 96971.83232729399, 96971.83232729399, 96971.83232729399)
$data5[9] = @(-0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02)
$data5[10] = @(-0.05739333478772612, -0.05930365251307265, -0.06092148188669505, -0.06232717353101374, -0.06356959090477614, -0.06468536632279165, -0.06569272067195379, -0.06661289774811661, -0.06745966700262382, -0.06824585120875767, -0.06871185127715168, -0.06873112831137719, -0.06873112831137719, -0.06873112831137719, -0.06873112831137719)

# --- Scenario data for RESUMEN E3 (rows 3-10, columns B-P) ---
$data6 = @{}
$data6[3] = @(708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708)
$data6[4] = @(14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674)
$data6[5] = @(43.46256804826548, 44.93474345481792, 46.19319584025903, 47.29463128040483, 48.27392091540726, 49.1577701117355, 49.95913435150531, 50.69385846266223, 51.37217991645741, 52.00381018296633, 52.27767271923103, 52.27767271923103, 52.27767271923103, 52.27767271923103, 52.27767271923103)
$data6[6] = @(29.01358845642875, 30.48576386298118, 31.74421624842229, 32.8456516885681, 33.82494132357053, 34.70879051989877, 35.51015475966858, 36.2448788708255, 36.92320032462068, 37.55483059112959, 37.8286931273943, 37.8286931273943, 37.8286931273943, 37.8286931273943, 37.8286931273943)
$data6[7] = @(0.04097964471246998, 0.04305898850703557, 0.0448364636277151, 0.04639216340193234, 0.04777534085250074, 0.04902371542358584, 0.05015558581874093, 0.0511933317384541, 0.05215141288788232, 0.05304354603266891, 0.05343035752456823, 0.05343035752456823, 0.05343035752456823, 0.05343035752456823, 0.05343035752456823)
$data6[8] = @(51057.79394630823, 53648.51203943407, 55863.12271000264, 57801.41668672254, 59524.75981860704, 61080.14791592934, 62490.37989366462, 63783.33930580939, 64977.04194170518, 66088.57793960819, 66570.51821966752, 66570.51821966752, 66570.51821966752, 66570.51821966752, 66570.51821966752)
$data6[9] = @(-0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02)
$data6[10] = @(-0.05783730274303422, -0.0596794660432805, -0.0612484918917817, -0.06261745989141897, -0.06383126481073902, -0.06492407798242786, -0.06591270173721087, -0.06681727800641878, -0.06765085853172702, -0.06842572298479231, -0.06876128892783764, -0.06876128892783764, -0.06876128892783764, -0.06876128892783764, -0.06876128892783764)

# --- Scenario data for RESUMEN E4 (rows 3-10, columns B-P) ---
$data7 = @{}
$data7[3] = @(708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708, 708)
$data7[4] = @(14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674, 14.44897959183674)
$data7[5] = @(43.46256804826548, 44.93474345481792, 46.19319584025903, 47.29463128040483, 48.27392091540726, 49.1577701117355, 49.95913435150531, 50.69385846266223, 51.37217991645741, 52.00381018296633, 52.27767271923103, 52.27767271923103, 52.27767271923103, 52.27767271923103, 52.27767271923103)
$data7[6] = @(29.01358845642875, 30.48576386298118, 31.74421624842229, 32.8456516885681, 33.82494132357053, 34.70879051989877, 35.51015475966858, 36.2448788708255, 36.92320032462068, 37.55483059112959, 37.8286931273943, 37.8286931273943, 37.8286931273943, 37.8286931273943, 37.8286931273943)
$data7[7] = @(0.04097964471246998, 0.04305898850703557, 0.0448364636277151, 0.04639216340193234, 0.04777534085250074, 0.04902371542358584, 0.05015558581874093, 0.0511933317384541, 0.05215141288788232, 0.05304354603266891, 0.05343035752456823, 0.05343035752456823, 0.05343035752456823, 0.05343035752456823, 0.05343035752456823)
$data7[8] = @(51057.79394630823, 53648.51203943407, 55863.12271000264, 57801.41668672254, 59524.75981860704, 61080.14791592934, 62490.37989366462, 63783.33930580939, 64977.04194170518, 66088.57793960819, 66570.51821966752, 66570.51821966752, 66570.51821966752, 66570.51821966752, 66570.51821966752)
$data7[9] = @(-0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02, -0.02)
$data7[10] = @(-0.05783730274303422, -0.0596794660432805, -0.0612484918917817, -0.06261745989141897, -0.06383126481073902, -0.06492407798242786, -0.06591270173721087, -0.06681727800641878, -0.06765085853172702, -0.06842572298479231, -0.06876128892783764, -0.06876128892783764, -0.06876128892783764, -0.06876128892783764, -0.06876128892783764)

# --- Duplicate RESUMEN E1 (carrying over all formatting) to build the new ---
# --- scenario sheets, placed in order right after RESUMEN E1.            ---
$newSheetNames = @('RESUMEN E2', 'RESUMEN E3', 'RESUMEN E4')
$newSheetData = @($data5, $data6, $data7)
$scenarioLabels = @(
    "Escenario: medidores con CAPEX ≈ 65% del CAPEX de E1",
    "Escenario: medidores con CAPEX ≈ 15% del CAPEX de E1",
    "Escenario: medidores con CAPEX ≈ 300000",
)

$afterSheet = $ws1
for ($i = 0; $i -lt $newSheetNames.Count; $i++) {
    $ws1.Copy($null, $afterSheet)
    $ws = $wb.Worksheets.Item($afterSheet.Index + 1)
    $ws.Name = $newSheetNames[$i]

    $ws.Cells.Item(2, 1).Value2 = $scenarioLabels[$i]

    $rowData = $newSheetData[$i]
    foreach ($r in $rowData.Keys) {
        $col = 2
        foreach ($v in $rowData[$r]) {
            $ws.Cells.Item($r, $col).Value2 = $v
            $col = $col + 1
        }
    }

    $afterSheet = $ws
}

Write-Output ("Sheets now: " + $wb.Worksheets.Count)
foreach ($s in $wb.Worksheets) { Write-Output $s.Name }
